# Auto-generated edit script: updates H:N market-profit columns per commit diff
# "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1309.1
$ws.Range("I32").Value = 775
$ws.Range("J32").Value = 1665.1666
$ws.Range("K32").Value = 775
$ws.Range("L32").Value = 1665.1666
$ws.Range("M32").Value = -449
$ws.Range("N32").Value = -2317.1666
$ws.Range("H42").Value = 153.14285
$ws.Range("I42").Value = 125
$ws.Range("J42").Value = 168.77777
$ws.Range("K42").Value = 375
$ws.Range("L42").Value = 506.33331
$ws.Range("M42").Value = -145
$ws.Range("N42").Value = -966.33331
$ws.Range("H99").Value = 1311.2
$ws.Range("I99").Value = 203.33333
$ws.Range("J99").Value = 2049.7778
$ws.Range("K99").Value = 609.99999
$ws.Range("L99").Value = 6149.3334
$ws.Range("M99").Value = 888.00001
$ws.Range("N99").Value = -9145.3334
$ws.Range("H116").Value = 16144.889
$ws.Range("I116").Value = 29176.25
$ws.Range("K116").Value = 29176.25
$ws.Range("M116").Value = -25734.25
$ws.Range("H129").Value = 1025.6786
$ws.Range("J129").Value = 1051.4321
$ws.Range("L129").Value = 3154.2963
$ws.Range("N129").Value = -13154.2963

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 323685.03
$ws.Range("I32").Value = 4483.4414
$ws.Range("J32").Value = 2214340.5
$ws.Range("K32").Value = 4483.4414
$ws.Range("L32").Value = 2214340.5
$ws.Range("M32").Value = -4196.4414
$ws.Range("N32").Value = -2214914.5
$ws.Range("H37").Value = 8199.799999999999
$ws.Range("J37").Value = 8199.799999999999
$ws.Range("L37").Value = 8199.799999999999
$ws.Range("N37").Value = -8745.799999999999
$ws.Range("H61").Value = 3361.2327
$ws.Range("I61").Value = 3742.9375
$ws.Range("J61").Value = 2250.818
$ws.Range("K61").Value = 3742.9375
$ws.Range("L61").Value = 2250.818
$ws.Range("M61").Value = -3530.9375
$ws.Range("N61").Value = -2674.818
$ws.Range("H74").Value = 1283.6666
$ws.Range("I74").Value = 815.8461
$ws.Range("K74").Value = 815.8461
$ws.Range("M74").Value = 58.15390000000002
$ws.Range("H77").Value = 1283.6666
$ws.Range("I77").Value = 815.8461
$ws.Range("K77").Value = 4079.2305
$ws.Range("M77").Value = 288.7695000000003
$ws.Range("H102").Value = 3629.2144
$ws.Range("I102").Value = 2812.3333
$ws.Range("J102").Value = 5099.6
$ws.Range("K102").Value = 2812.3333
$ws.Range("L102").Value = 5099.6
$ws.Range("M102").Value = -1190.3333
$ws.Range("N102").Value = -8343.6
$ws.Range("H132").Value = 1579.1384
$ws.Range("I132").Value = 1273.5862
$ws.Range("J132").Value = 4110.857
$ws.Range("K132").Value = 3820.7586
$ws.Range("L132").Value = 12332.571
$ws.Range("M132").Value = -1290.7586
$ws.Range("N132").Value = -17392.571
$ws.Range("H133").Value = 78646
$ws.Range("J133").Value = 78646
$ws.Range("L133").Value = 78646
$ws.Range("N133").Value = -83706
$ws.Range("H136").Value = 3361.2327
$ws.Range("I136").Value = 3742.9375
$ws.Range("J136").Value = 2250.818
$ws.Range("K136").Value = 11228.8125
$ws.Range("L136").Value = 6752.454000000001
$ws.Range("M136").Value = -8678.8125
$ws.Range("N136").Value = -11852.454
$ws.Range("H138").Value = 62520
$ws.Range("J138").Value = 62520
$ws.Range("L138").Value = 62520
$ws.Range("N138").Value = -72800

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1694.1177
$ws.Range("I105").Value = 1693.3334
$ws.Range("K105").Value = 1693.3334
$ws.Range("M105").Value = 53.66660000000002
$ws.Range("H138").Value = 59633.332
$ws.Range("J138").Value = 59633.332
$ws.Range("L138").Value = 59633.332
$ws.Range("N138").Value = -69913.33199999999
$ws.Range("H140").Value = 89750
$ws.Range("J140").Value = 89750
$ws.Range("L140").Value = 89750
$ws.Range("N140").Value = -100110

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 8250.333000000001
$ws.Range("J60").Value = 8250.333000000001
$ws.Range("L60").Value = 8250.333000000001
$ws.Range("N60").Value = -9272.333000000001
$ws.Range("H68").Value = 17201
$ws.Range("J68").Value = 17201
$ws.Range("L68").Value = 17201
$ws.Range("N68").Value = -18699
$ws.Range("H71").Value = 17201
$ws.Range("J71").Value = 17201
$ws.Range("L71").Value = 51603
$ws.Range("N71").Value = -59091
$ws.Range("H74").Value = 15224.333
$ws.Range("J74").Value = 17812.2
$ws.Range("L74").Value = 17812.2
$ws.Range("N74").Value = -19560.2
$ws.Range("H77").Value = 15224.333
$ws.Range("J77").Value = 17812.2
$ws.Range("L77").Value = 53436.60000000001
$ws.Range("N77").Value = -62172.60000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 483.85715
$ws.Range("I107").Value = 258.33334
$ws.Range("J107").Value = 653
$ws.Range("K107").Value = 775.0000200000001
$ws.Range("L107").Value = 1959
$ws.Range("M107").Value = 1144.99998
$ws.Range("N107").Value = -5799
$ws.Range("H113").Value = 999.2174
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 994.75555
$ws.Range("K113").Value = 3600
$ws.Range("L113").Value = 2984.26665
$ws.Range("M113").Value = -1430
$ws.Range("N113").Value = -7324.26665
$ws.Range("H122").Value = 491.56668
$ws.Range("I122").Value = 365.92
$ws.Range("K122").Value = 3293.28
$ws.Range("M122").Value = -843.2800000000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 53190
$ws.Range("J133").Value = 53190
$ws.Range("L133").Value = 53190
$ws.Range("N133").Value = -63310
$ws.Range("H138").Value = 67914.28999999999
$ws.Range("J138").Value = 67914.28999999999
$ws.Range("L138").Value = 67914.28999999999
$ws.Range("N138").Value = -78194.28999999999
$ws.Range("H139").Value = 65666.664
$ws.Range("J139").Value = 65666.664
$ws.Range("L139").Value = 65666.664
$ws.Range("N139").Value = -75946.664
$ws.Range("H140").Value = 99894.5
$ws.Range("J140").Value = 99894.5
$ws.Range("L140").Value = 99894.5
$ws.Range("N140").Value = -110254.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2864.1052
$ws.Range("I122").Value = 2450.7
$ws.Range("J122").Value = 3323.4443
$ws.Range("K122").Value = 7352.099999999999
$ws.Range("L122").Value = 9970.332900000001
$ws.Range("M122").Value = -4902.099999999999
$ws.Range("N122").Value = -14870.3329
$ws.Range("H132").Value = 5062.222
$ws.Range("I132").Value = 6108.5835
$ws.Range("K132").Value = 18325.7505
$ws.Range("M132").Value = -15795.7505
$ws.Range("H139").Value = 70400
$ws.Range("J139").Value = 70400
$ws.Range("L139").Value = 70400
$ws.Range("N139").Value = -80680

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1026.6
$ws.Range("I100").Value = 783.25
$ws.Range("K100").Value = 1566.5
$ws.Range("M100").Value = -1025.5
$ws.Range("H136").Value = 998.6
$ws.Range("I136").Value = 909.1852
$ws.Range("J136").Value = 1803.3334
$ws.Range("K136").Value = 2727.5556
$ws.Range("L136").Value = 5410.0002
$ws.Range("M136").Value = -177.5556000000001
$ws.Range("N136").Value = -10510.0002

